# Fruta / hortaliza, semanal
# Insert a new weekly record for "Vega Modelo de Temuco - Chirimoya" above the
# existing row 66, pushing the previously-recorded rows (old 66-98) down by
# one row (new 67-99), and fill the newly-inserted row 66 with this week's
# price data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 66:98 down to 67:99, inserting a fresh blank row at 66.
$ws.Rows.Item(66).Insert()

# Populate the new row 66 with the latest weekly Chirimoya price record.
$ws.Range("A66").Value = 10
$ws.Range("B66").Value = "Vega Modelo de Temuco"
$ws.Range("C66").Value = "La Araucanía"
$ws.Range("D66").Value = 44510
$ws.Range("E66").Value = 9
$ws.Range("F66").Value = "Fruta"
$ws.Range("G66").Value = 100107
$ws.Range("H66").Value = "Otros"
$ws.Range("I66").Value = 100107002
$ws.Range("J66").Value = "Chirimoya"
$ws.Range("K66").Value = "Cultivar IV Región"
$ws.Range("L66").Value = "Primera"
$ws.Range("M66").Value = 35
$ws.Range("N66").Value = 3000
$ws.Range("O66").Value = 3000
$ws.Range("P66").Value = 3000
$ws.Range("Q66").Value = "$/kilo (en caja de 15 kilos)"
$ws.Range("R66").Value = "Provincia del Elquí"
$ws.Range("S66").Value = 3000
$ws.Range("T66").Value = 1
